$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (rows 2-5) from numeric placeholder values to text labels A/B/C/A
$ws.Range("F2").Value = "A"
$ws.Range("F3").Value = "B"
$ws.Range("F4").Value = "C"
$ws.Range("F5").Value = "A"

# Move the active selection from C3 to F2
$ws.Range("F2").Select()
